# Auto-generated: updates cryptos price/volume table cell values to match
# the target snapshot (GitHub Actions crypto-price refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.409.18"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").Value = "1.837.11"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'230.60"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'43.42"
$ws.Range("E8").Value = "  +13.42%  "
$ws.Range("E9").Value = "  +6.95%  "
$ws.Range("D10").Value = "'0.0699"
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("E11").Value = "  +2.93%  "
$ws.Range("D12").Value = "2.101.89"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "1.835.07"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "'11.24"
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "'0.670"
$ws.Range("E15").Value = "  +6.42%  "
$ws.Range("E16").Value = "  +6.59%  "
$ws.Range("D17").Value = "35.355.56"
$ws.Range("E17").Value = "  +2.62%  "
$ws.Range("D18").Value = "'69.93"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("D20").Value = "'244.60"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").Value = "'12.02"
$ws.Range("E21").Value = "  +7.91%  "
$ws.Range("D22").Value = "'4.67"
$ws.Range("E22").Value = "  +13.89%  "
$ws.Range("D23").Value = "'1.01"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "'169.19"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").Value = "'7.88"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("D27").Value = "'17.70"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "'0.121"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "'1.52"
$ws.Range("E29").Value = "  +23.69%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "3.304.08"
$ws.Range("E31").Value = "  +35.99%  "
$ws.Range("D32").Value = "'0.0549"
$ws.Range("E32").Value = "  +6.87%  "
$ws.Range("E33").Value = "  +3.99%  "
$ws.Range("D34").Value = "'4.06"
$ws.Range("E34").Value = "  +5.76%  "
$ws.Range("D35").Value = "'1.85"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").Value = "'95.73"
$ws.Range("E36").Value = "  +16.34%  "
$ws.Range("D37").Value = "'0.683"
$ws.Range("E37").Value = "  +6.68%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.09"
$ws.Range("E38").Value = "  +3.31%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.344.23"
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'15.33"
$ws.Range("E40").Value = "  +9.66%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0194"
$ws.Range("E41").Value = "  +4.16%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.42"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = "  +6.06%  "
$ws.Range("E44").Value = "  +3.71%  "
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'6.21"
$ws.Range("E47").Value = "  +7.24%  "
$ws.Range("D48").Value = "'0.0520"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "2.004.16"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "'102.50"
$ws.Range("E51").Value = "  -0.24%  "
